# Fonds de solidarite - add 2020-07-16 data (Mayotte) + refresh cumulative
# counters for several regions already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh nombre_aides (C) / montant_total (D) for rows that were
#    already in the sheet before the new Mayotte rows are inserted.
#    (Values stay text, like the rest of the sheet, via the leading
#    quote so Excel doesn't coerce them to numbers.)
# ---------------------------------------------------------------------
function Set-Counts($row, $count, $amount) {
    $ws.Cells.Item($row, 3).Value = "'" + $count
    $ws.Cells.Item($row, 4).Value = "'" + $amount
}

Set-Counts 2  "388" "880328.79"
Set-Counts 4  "814" "2397771.47"
Set-Counts 6  "520" "1436121.69"
Set-Counts 9  "114" "253455.66"
Set-Counts 13 "166" "391083.00"
Set-Counts 15 "410" "1296754.00"
Set-Counts 17 "10"  "22850.00"
Set-Counts 19 "134" "312000.00"
Set-Counts 20 "285" "847336.33"
Set-Counts 21 "140" "352900.26"
Set-Counts 27 "184" "408926.00"
Set-Counts 28 "4"   "17500.00"
Set-Counts 29 "385" "1161532.00"
Set-Counts 31 "299" "821218.74"
Set-Counts 42 "190" "482472.74"
Set-Counts 44 "410" "1207333.16"
Set-Counts 45 "278" "755762.76"
Set-Counts 47 "17"  "62220.65"

# ---------------------------------------------------------------------
# 2) Insert two new rows right before the old row 59 ("Normandie") to
#    hold the new "Mayotte" (reg 06) entries, pushing every row from
#    the old 59 onward down by two.
# ---------------------------------------------------------------------
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

# ---------------------------------------------------------------------
# 3) Populate the two new Mayotte rows.
# ---------------------------------------------------------------------
function Set-Row($row, $dispositif, $volet, $count, $amount, $reg, $libelleRegion, $catCode, $catLibelle) {
    $ws.Cells.Item($row, 1).Value = $dispositif
    $ws.Cells.Item($row, 2).Value = $volet
    $ws.Cells.Item($row, 3).Value = "'" + $count
    $ws.Cells.Item($row, 4).Value = "'" + $amount
    $ws.Cells.Item($row, 5).Value = "'" + $reg
    $ws.Cells.Item($row, 6).Value = $libelleRegion
    $ws.Cells.Item($row, 7).Value = "'" + $catCode
    $ws.Cells.Item($row, 8).Value = $catLibelle
}

Set-Row 59 "Fonds de solidarité" "VOLET2" "3" "7500.00" "06" "Mayotte" "10" "Entrepreneur individuel"
Set-Row 60 "Fonds de solidarité" "VOLET2" "4" "8000.00" "06" "Mayotte" "54" "Société à responsabilité limitée (SARL)"

# ---------------------------------------------------------------------
# 4) Refresh nombre_aides (C) / montant_total (D) for the
#    Nouvelle-Aquitaine rows, which now sit two rows further down.
# ---------------------------------------------------------------------
Set-Counts 66 "304" "703118.84"
Set-Counts 68 "758" "2256051.03"
Set-Counts 69 "432" "1196267.79"
Set-Counts 70 "31"  "83000.00"
Set-Counts 71 "26"  "78231.09"
